# Applies the edit described by the diff:
#   - The two rows that describe "생년/age" (birth year) and
#     "학번/hakbun" (student id) swap places: what used to be row 7
#     (생년/age/Int 4자릿수/1996~2003) becomes row 8, and what used to be
#     row 8 (학번/hakbun/Int 그대로/14~21) becomes row 7. The "mmm-yy"
#     number format that was on E7 stays with the E column content that
#     ends up in row 8 (it tracks the row position, not the text), so the
#     format must travel along with the swap.
#   - The sheet's selection moves to the edited row (A8:E8).
#   - Page setup (paper size / orientation) gets recorded, matching the
#     <pageSetup paperSize="9" orientation="portrait".../> added to the
#     sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read current row 7 / row 8 contents (Value2 is reliable for reads;
#     the plain .Value getter is not dependable for readback in this
#     runtime) ---
$a7 = $ws.Range("A7").Value2
$c7 = $ws.Range("C7").Value2
$d7 = $ws.Range("D7").Value2
$e7 = $ws.Range("E7").Value2

$a8 = $ws.Range("A8").Value2
$c8 = $ws.Range("C8").Value2
$d8 = $ws.Range("D8").Value2
$e8 = $ws.Range("E8").Value2

# --- Write row 8's old content into row 7 ---
$ws.Range("A7").Value = $a8
$ws.Range("C7").Value = $c8
$ws.Range("D7").Value = $d8
$ws.Range("E7").Value = $e8

# --- Write row 7's old content into row 8 ---
$ws.Range("A8").Value = $a7
$ws.Range("C8").Value = $c7
$ws.Range("D8").Value = $d7
$ws.Range("E8").Value = $e7

# --- The date-ish "mmm-yy" number format stays on row 8 / column E now,
#     and row 7's E cell goes back to the default (General) format ---
$ws.Range("E8").NumberFormat = "mmm-yy"
$ws.Range("E7").ClearFormats()

# --- Selection now highlights the newly-written row 8 ---
[void]$ws.Range("A8:E8").Select()

# --- Page setup touched: letter/A4-ish paper (id 9 = A4), portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
